# "add reward and agent final adjustments"
#
# hyperparameters sheet: bump a couple of training hyperparameters and
# move the cell selection to where the user left off (C4).
# generator_parameters sheet: scroll the view down a bit (top-left cell A4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("hyperparameters")
$ws2 = $wb.Worksheets.Item("generator_parameters")

# n_updates_per_iteration: 1 -> 12
$ws1.Range("B5").Value = 12

# iteration_number: 100 -> 10000
$ws1.Range("B11").Value = 10000

# Scroll the generator_parameters view so row 4 is at the top, leaving the
# existing selection (B17) untouched.
$ws2.Activate()
$ws2.Application.ActiveWindow.ScrollRow = 4
$ws2.Application.ActiveWindow.ScrollColumn = 1

# Leave the hyperparameters sheet active/selected, with C4 as the
# last-selected cell, matching where the author finished editing.
$ws1.Activate()
$ws1.Range("C4").Select()
